$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class")

# Populate the text cells in the same order the original workbook's
# shared-string table was built (first-seen order: Status, U49,
# TestAlluree, active, Classno, " Classdate", Batchname, Staffname,
# Classtopic), so the resulting shared-string indices line up too.
$ws.Range("F1").Value = "Status"
$ws.Range("B2").Value = "U49"
$ws.Range("D2").Value = "TestAlluree"
$ws.Range("F2").Value = "active"
$ws.Range("C1").Value = "Classno"
$ws.Range("E1").Value = " Classdate"
$ws.Range("A1").Value = "Batchname"
$ws.Range("B1").Value = "Staffname"
$ws.Range("D1").Value = "Classtopic"

# Numeric cells
$ws.Range("A2").Value = 8547
$ws.Range("C2").Value = 1

# Set the number format before assigning the date value so the engine
# registers the built-in date format (numFmtId 14) instead of inventing
# a custom one.
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Value = [datetime]"2024-12-23"

# Column widths (approximate autofit result from the original edit; the
# engine snaps ColumnWidth to 1/6-character increments, so these are the
# closest inputs to the recorded widths of 11.855.., 10.711.., 12.711..,
# 11.285.. characters).
$ws.Columns.Item(1).ColumnWidth = 11.0
$ws.Columns.Item(2).ColumnWidth = 9.833333333333334
$ws.Columns.Item(4).ColumnWidth = 11.833333333333334
$ws.Columns.Item(5).ColumnWidth = 10.5

$ws.Range("E9").Select()
